$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.649.61'
$ws.Range('D2').Style = 'Normal'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.598.03'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.17%  '
$ws.Range('E4').Value = '  +0.12%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '211.28'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.14%  '
$ws.Range('E6').Value = '  +0.82%  '
$ws.Range('E7').Value = '  +0.07%  '
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('E9').Value = '  +0.23%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.54'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.55%  '
$ws.Range('E11').Value = '  +0.62%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.822.07'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.22%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.611.90'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.38%  '
$ws.Range('E14').Value = '  +0.18%  '
$ws.Range('E15').Value = '  +0.17%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '64.85'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.02%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '26.630.51'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.02%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.0₃0735'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.64%  '
$ws.Range('E19').Value = '  +0.07%  '
$ws.Range('E20').Value = '  -0.46%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '7.05'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +5.48%  '
$ws.Range('E22').Value = '  +0.75%  '
$ws.Range('E23').Value = '  +0.58%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '8.94'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.40%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '145.50'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.71%  '
$ws.Range('E26').Value = '  +0.07%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.12'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.16%  '
$ws.Range('E28').Value = '  +0.52%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0513'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +2.13%  '
$ws.Range('E31').Value = '  +0.17%  '
$ws.Range('E32').Value = '  +0.27%  '
$ws.Range('E33').Value = '  +1.04%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.274.98'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.61%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.621'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -10.13%  '
$ws.Range('E36').Value = '  +0.54%  '
$ws.Range('E37').Value = '  +0.67%  '
$ws.Range('E38').Value = '  -0.41%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.08'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +19.27%  '
$ws.Range('E40').Value = '  +0.07%  '
$ws.Range('E41').Value = '  +2.95%  '
$ws.Range('E42').Value = '  +0.54%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.785'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.87%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '63.89'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.60%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.734.66'
$ws.Range('D45').Style = 'Normal'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '90.18'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E47').Value = '  -2.30%  '
$ws.Range('E48').Value = '  +3.25%  '
$ws.Range('E49').Value = '  +1.00%  '
$ws.Range('E50').Value = '  +0.01%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.44'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.47%  '
